$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the three new "ara" (Arabic) title rows (8, 9, 10) under the existing
# eng / fra rows, mirroring the Mr / Mrs / Miss pattern already present.
# ---------------------------------------------------------------------------

# Row 8: ara / MIR / "Mr" / "Male Title"
$ws.Range("A8").Value = "ara"
$ws.Range("B8").Value = "MIR"
$ws.Range("C8").Value = "السيد"
$ws.Range("D8").Value = "لقب ذكر"

# Row 9: ara / MRS / "Mrs" / "Female Title"
$ws.Range("A9").Value = "ara"
$ws.Range("B9").Value = "MRS"
$ws.Range("C9").Value = "السيدة"
$ws.Range("D9").Value = "لقب أنثى"

# Row 10: ara / MIS / "Miss" / "Unmarried Female Title"
$ws.Range("A10").Value = "ara"
$ws.Range("B10").Value = "MIS"
$ws.Range("C10").Value = "يغيب"
$ws.Range("D10").Value = "لقب انثى غير متزوجة"

# ---------------------------------------------------------------------------
# Column E ("is_active") must stay a text "TRUE" (shared string), not a
# boolean. Copy the existing E7 cell (already text "TRUE") down instead of
# assigning the literal value, which Excel would otherwise coerce to a bool.
# ---------------------------------------------------------------------------
$ws.Range("E7").Copy($ws.Range("E8"))
$ws.Range("E7").Copy($ws.Range("E9"))
$ws.Range("E7").Copy($ws.Range("E10"))

# ---------------------------------------------------------------------------
# Columns C and D in the new rows get a left aligned, word-wrapped style.
# Build the style once on C8 (wrap first, then alignment, keeps the
# stylesheet from growing extra transient entries), then propagate that
# same format to the rest of the C:D block via a formats-only paste so the
# other cells reuse the same style instead of minting new ones.
# ---------------------------------------------------------------------------
$ws.Range("C8").WrapText = $true
$ws.Range("C8").HorizontalAlignment = -4131

$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("D10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Widen columns C, D, E to fit the new Arabic text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 22.766666666666667
$ws.Columns.Item(4).ColumnWidth = 24.636666666666667
$ws.Columns.Item(5).ColumnWidth = 17.356666666666667

# ---------------------------------------------------------------------------
# Taller rows (wrapped Arabic text) for the new data rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 16.4
$ws.Rows.Item(9).RowHeight = 16.4
$ws.Rows.Item(10).RowHeight = 16.4

# ---------------------------------------------------------------------------
# Leave the selection where the author's session ended up.
# ---------------------------------------------------------------------------
$ws.Range("I5").Select()
